# SKA-100: remove the wheels (WH-01) line item from the root assembly
# parts list — delete the entire row 3 and shift the rows below it up.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(3).Delete()
